$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 20 and 21 currently hold unfilled placeholder cells (blank, with the
# "white" row style used by rows 17-19). Bring them in line with the rest of
# the finished rows by copying the formatting from row 12 (one of the
# "blue" finished rows) and then filling in the new activity data.
$ws.Range("A12:F12").Copy()
$ws.Range("A20:F21").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Column C on these rows needs the date number format (built-in format 16,
# "d-mmm") instead of the plain formatting it inherited from row 12.
$ws.Range("C20:C21").NumberFormat = "d-mmm"

# Row 20: "Subir el programa a GITHUB"
$ws.Range("B20").Value = "Subir el programa a GITHUB"
$ws.Range("C20").Value = "20 de Febrero"
$ws.Range("D20").Value = "20 de Febrero"
$ws.Range("E20").Value = "Nayre"
$ws.Range("F20").Value = "Subido Finalizado"

# Row 21: "Subir Cronograma"
$ws.Range("B21").Value = "Subir Cronograma"
$ws.Range("C21").Value = "20 de Febrero"
$ws.Range("D21").Value = "20 de Febrero"
$ws.Range("E21").Value = "Nayre"
$ws.Range("F21").Value = "Subido Finalizado"

# Match the author's final on-screen selection (row 21 highlighted).
$null = $ws.Range("A21:F21").Select()
